$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Bills (SBM)" column (F) and the "IGs (SBM)" column (H).
# Deleting column H first keeps column F's index valid for the second delete.
$ws.Range("H1").EntireColumn.Delete()
$ws.Range("F1").EntireColumn.Delete()
